# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45178 (2023-09-09) to 45179 (2023-09-10), leaving everything else
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 203 }

$ws.Range("C2:C$lastRow").Value = 45179
